$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, pushing the existing price history rows down
# by one (row 2 -> 3, row 3 -> 4, ... row 18 -> 19).
$ws.Rows("2:2").Insert()

# Excel's row insert copies formatting from the row above (the bold header
# row), so strip that back off to match the plain formatting used by the
# rest of the data rows.
$ws.Range("A2:D2").ClearFormats()

# Populate the newly inserted row with today's price data. The date column
# is forced to text first so "2025-12-08" is stored as a literal string
# (like all the other date cells) instead of being auto-converted into a
# date serial number, then the style is reset to Normal so no stray
# text-format style lingers on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-08"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
